$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Disaster" column header in I1, matching style of H1 (font color black)
$ws.Range("I1").Value = "Disaster"
$ws.Range("I1").Font.Color = 0

# Fill in the Disaster column values (I2:I5), matching style of H column
$ws.Range("I2").Value = 0.25
$ws.Range("I3").Value = 0.05
$ws.Range("I4").Value = 0.1
$ws.Range("I5").Value = 0.075

$ws.Range("I2:I5").Font.Color = 0

# Update selection to I5 to match final cursor position
$ws.Range("I5").Select()
